$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A7").Value = 42602.582546296297
$ws.Range("A7").NumberFormat = "m/d/yy h:mm"

$ws.Range("B7").Value = "Random"

$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 51
$ws.Range("I7").Value = 49
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 70
$ws.Range("M7").Value = 30
